# Populate the previously-empty column A ("http_id" style row labels h1, h2, ...)
# for rows 21 through 151 of Sheet1 with values h20, h21, ..., h150, matching the
# pattern already used by the existing rows 2-20 (h1 .. h19).
# Column A cells elsewhere use style index 3 (Arial 10, left-aligned), so apply the
# equivalent formatting (left horizontal alignment) to the newly written cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 21; $row -le 151; $row++) {
    $label = "h" + ($row - 1)
    $cell = $ws.Cells.Item($row, 1)
    $cell.Value = $label
    $cell.HorizontalAlignment = -4131   # xlLeft
}

# Update the saved view state to match: scrolled so row 175 is at the top and
# cell C192 is selected.
$win = $excel.ActiveWindow
$win.ScrollRow = 175
$win.ScrollColumn = 1
$null = $ws.Range("C192").Select()
